$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 16:20:15"
$ws.Range("A3").Value = "Total filas: 314"
$ws.Range("A142").Value = "11:17:39"
$ws.Range("C142").Value = "225_C ROCA-H SUR"
$ws.Range("D142").Value = 3
$ws.Range("A143").Value = "10:13:53"
$ws.Range("C143").Value = "26_HERNANDEZ"
$ws.Range("D143").Value = 67
$ws.Range("A162").Value = "10:13:53"
$ws.Range("C162").Value = "16_P MOR-SANTA ANA"
$ws.Range("D162").Value = 113
$ws.Range("A163").Value = "10:52:37"
$ws.Range("C163").Value = "10_OLMOS"
$ws.Range("D163").Value = 74
$ws.Range("C164").Value = "14_ABASTO"
$ws.Range("A173").Value = "12:01:11"
$ws.Range("C173").Value = "14_ABASTO"
$ws.Range("D173").Value = 20
$ws.Range("A174").Value = "10:52:37"
$ws.Range("C174").Value = "26_HERNANDEZ"
$ws.Range("D174").Value = 89
$ws.Range("C175").Value = "215A_EL PATO"
$ws.Range("A193").Value = "10:52:37"
$ws.Range("C193").Value = "15_ABASTO"
$ws.Range("D193").Value = 118
$ws.Range("A194").Value = "12:50:41"
$ws.Range("C194").Value = "16_SANTA ANA"
$ws.Range("D194").Value = 0
$ws.Range("A202").Value = "12:50:41"
$ws.Range("C202").Value = "11_ETCHEVERRY"
$ws.Range("D202").Value = 24
$ws.Range("A203").Value = "12:01:11"
$ws.Range("C203").Value = "215D_EL PATO"
$ws.Range("D203").Value = 73
$ws.Range("A220").Value = "12:01:11"
$ws.Range("C220").Value = "215A_EL PATO"
$ws.Range("D220").Value = 110
$ws.Range("A221").Value = "13:51:32"
$ws.Range("C221").Value = "11_ETCHEVERRY"
$ws.Range("D221").Value = 0
$ws.Range("C237").Value = "10_OLMOS"
$ws.Range("C238").Value = "16_SANTA ANA"
$ws.Range("C239").Value = "11_ETCHEVERRY"
$ws.Range("C268").Value = "16_P MOR-167 Y 521"
$ws.Range("C269").Value = "14_ABASTO"
$ws.Range("A292").Value = "16:20:15"
$ws.Range("B292").Value = "16:24"
$ws.Range("C292").Value = "14_ABASTO"
$ws.Range("D292").Value = 4
$ws.Range("A293").Value = "15:36:13"
$ws.Range("B293").Value = "16:29"
$ws.Range("C293").Value = "10_OLMOS"
$ws.Range("D293").Value = 53
$ws.Range("A294").Value = "14:49:07"
$ws.Range("B294").Value = "16:30"
$ws.Range("C294").Value = "15_ABASTO"
$ws.Range("D294").Value = 101
$ws.Range("A295").Value = "15:59:02"
$ws.Range("B295").Value = "16:35"
$ws.Range("C295").Value = "23_HERNANDEZ"
$ws.Range("D295").Value = 36
$ws.Range("A296").Value = "16:20:15"
$ws.Range("B296").Value = "16:35"
$ws.Range("C296").Value = "16_SANTA ANA"
$ws.Range("D296").Value = 15
$ws.Range("A297").Value = "15:36:13"
$ws.Range("B297").Value = "16:37"
$ws.Range("C297").Value = "11_ETCHEVERRY"
$ws.Range("D297").Value = 61
$ws.Range("A298").Value = "15:36:13"
$ws.Range("B298").Value = "16:40"
$ws.Range("C298").Value = "17_ROMERO"
$ws.Range("D298").Value = 64
$ws.Range("A299").Value = "14:59:23"
$ws.Range("B299").Value = "16:42"
$ws.Range("C299").Value = "16_P MOR-SANTA ANA"
$ws.Range("D299").Value = 103
$ws.Range("C300").Value = "16_P MOR-SANTA ANA"
$ws.Range("B301").Value = "16:43"
$ws.Range("C301").Value = "23_HERNANDEZ"
$ws.Range("D301").Value = 67
$ws.Range("A302").Value = "14:49:07"
$ws.Range("B302").Value = "16:43"
$ws.Range("C302").Value = "225_GOMEZ"
$ws.Range("D302").Value = 114
$ws.Range("A303").Value = "15:36:13"
$ws.Range("B303").Value = "16:48"
$ws.Range("C303").Value = "15_ABASTO"
$ws.Range("D303").Value = 72
$ws.Range("B304").Value = "16:51"
$ws.Range("C304").Value = "14_ABASTO"
$ws.Range("D304").Value = 52
$ws.Range("A305").Value = "14:59:23"
$ws.Range("B305").Value = "16:56"
$ws.Range("C305").Value = "17_179 Y 38"
$ws.Range("D305").Value = 117
$ws.Range("B306").Value = "16:57"
$ws.Range("C306").Value = "10_OLMOS"
$ws.Range("D306").Value = 58
$ws.Range("B307").Value = "17:05"
$ws.Range("C307").Value = "215A_EL PATO"
$ws.Range("D307").Value = 89
$ws.Range("A308").Value = "16:20:15"
$ws.Range("B308").Value = "17:05"
$ws.Range("C308").Value = "23_HERNANDEZ"
$ws.Range("D308").Value = 45
$ws.Range("B309").Value = "17:17"
$ws.Range("C309").Value = "11_ETCHEVERRY"
$ws.Range("D309").Value = 78
$ws.Range("A310").Value = "15:36:13"
$ws.Range("B310").Value = "17:21"
$ws.Range("C310").Value = "26_HERNANDEZ"
$ws.Range("D310").Value = 105
$ws.Range("A311").Value = "16:20:15"
$ws.Range("B311").Value = "17:21"
$ws.Range("C311").Value = "16_SANTA ANA"
$ws.Range("D311").Value = 61
$ws.Range("A312").Value = "15:36:13"
$ws.Range("B312").Value = "17:24"
$ws.Range("C312").Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Range("D312").Value = 108
$ws.Range("A313").Value = "16:20:15"
$ws.Range("B313").Value = "17:36"
$ws.Range("C313").Value = "27_EL RETIRO"
$ws.Range("D313").Value = 76
$ws.Range("A314").Value = "15:59:02"
$ws.Range("B314").Value = "17:37"
$ws.Range("C314").Value = "27_EL RETIRO"
$ws.Range("D314").Value = 98
$ws.Range("E314").Value = "LP1912"
$ws.Range("A315").Value = "15:59:02"
$ws.Range("B315").Value = "17:38"
$ws.Range("C315").Value = "17_ROMERO"
$ws.Range("D315").Value = 99
$ws.Range("E315").Value = "LP1912"
$ws.Range("A316").Value = "15:59:02"
$ws.Range("B316").Value = "17:40"
$ws.Range("C316").Value = "215B_EL PATO"
$ws.Range("D316").Value = 101
$ws.Range("E316").Value = "LP1912"
$ws.Range("A317").Value = "15:59:02"
$ws.Range("B317").Value = "17:51"
$ws.Range("C317").Value = "16_P MOR-167 Y 521"
$ws.Range("D317").Value = 112
$ws.Range("E317").Value = "LP1912"
$ws.Range("A318").Value = "15:59:02"
$ws.Range("B318").Value = "17:52"
$ws.Range("C318").Value = "81_EL PELIGRO"
$ws.Range("D318").Value = 113
$ws.Range("E318").Value = "LP1912"
$ws.Range("A319").Value = "16:20:15"
$ws.Range("B319").Value = "18:04"
$ws.Range("C319").Value = "17_ROMERO"
$ws.Range("D319").Value = 104
$ws.Range("E319").Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 16:20:15"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 16:20:15"
$ws.Range("A3").Value = "Total filas: 45"
$ws.Range("A49").Value = "16:20:15"
$ws.Range("B49").Value = "17:16"
$ws.Range("C49").Value = "215A_LA PLATA"
$ws.Range("D49").Value = 56
$ws.Range("E49").Value = "L6173"
$ws.Range("A50").Value = "16:20:15"
$ws.Range("B50").Value = "18:04"
$ws.Range("C50").Value = "215C_LA PLATA"
$ws.Range("D50").Value = 104
$ws.Range("E50").Value = "L6203"
